$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.692.66'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.55%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.742.94'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -5.51%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.98'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -8.70%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.05%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4930'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -6.47%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.58'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -7.83%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2421'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -23.30%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.05987'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -12.07%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.742.76'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -5.51%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06841'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -12.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.79'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -22.60%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.460'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -11.12%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.15'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -12.74%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.5811'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -26.06%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.01%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9998'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.05%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '25.733.95'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.52%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.50'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -17.42%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006458'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -18.58%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.962.92'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -5.64%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.973'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -13.89%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.012'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -16.55%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -16.12%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '136.21'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.86%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.471'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -12.58%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.845'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -17.06%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '14.56'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -14.61%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '101.00'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -9.04%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.773'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -10.59%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08101'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -6.94%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.353'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -18.03%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04391'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -10.16%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9987'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.07%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.633'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -8.05%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.019'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -10.72%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6075'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -16.94%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -13.27%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -12.32%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9994'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.08%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '102.83'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -6.51%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.01496'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -14.04%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7745'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -14.96%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.136'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -13.39%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -21.98%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05116'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -12.25%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.989'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -22.47%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -14.02%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '30.14'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -13.67%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.76'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -12.15%  '
